$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

$ws.Range("C2").Value = "CA"
$ws.Range("D2").Value = "Ontario"
$ws.Range("E2").Value = "Toronto"
$ws.Range("F2").Value = "Marvel pro, Finch Ave W, M3J 3H7, +14166041496"

$ws.Range("F2").Select()
